$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A17 is a date-looking string ("05/08/2025") that must stay literal text,
# not get auto-converted into a date serial number. Force text format for
# the assignment, then restore the cell's style so no stray formatting is
# left behind (matches the source row's unstyled cells).
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "05/08/2025"
$ws.Range("A17").Style = "Normal"

$ws.Range("B17").Value = "Santos"
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "Juventude"
$ws.Range("F17").Value = "L"
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 1.63
$ws.Range("L17").Value = 2.63
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 19
$ws.Range("O17").Value = 10
$ws.Range("P17").Value = 8
